$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New standings data for Rodada 6 (rows 2-21), columns A-H:
# A=Clube ID, B=Jogos, C=Pontos, D=Gols Marcados, E=Gols Sofridos, F=SG Total, G=Peso Casa, H=Peso Fora
$data = @(
    @("FLA", 6, 14, 15, 2, 13, 178.6, 137.2),
    @("PAL", 6, 13, 7, 3, 4, 160, 129.4),
    @("RBB", 6, 13, 8, 5, 3, 139, 109.8),
    @("FLU", 6, 10, 6, 6, 0, 138.4, 105.6),
    @("CRU", 6, 10, 7, 6, 1, 119.2, 85.40000000000001),
    @("INT", 6, 9, 8, 4, 4, 152.6, 112.8),
    @("BAH", 6, 9, 6, 7, -1, 138, 104.4),
    @("SAO", 6, 8, 6, 5, 1, 126.4, 96.2),
    @("CEA", 6, 8, 8, 7, 1, 106.8, 75.8),
    @("BOT", 6, 8, 6, 4, 2, 131.2, 92),
    @("MIR", 6, 7, 11, 9, 2, 105, 77.2),
    @("COR", 6, 7, 6, 10, -4, 128.6, 90),
    @("VAS", 6, 7, 6, 8, -2, 108.2, 71.59999999999999),
    @("JUV", 6, 7, 7, 14, -7, 107.2, 69.59999999999999),
    @("CAM", 6, 6, 6, 8, -2, 117.8, 81),
    @("FOR", 6, 6, 5, 5, 0, 105.2, 74.8),
    @("VIT", 6, 6, 7, 9, -2, 94.59999999999999, 64.19999999999999),
    @("GRE", 6, 5, 5, 11, -6, 81.59999999999999, 53.8),
    @("SAN", 6, 4, 7, 9, -2, 100.8, 66),
    @("SPT", 6, 2, 3, 8, -5, 78.59999999999999, 54.4)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $row++
}
